$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6136
$ws.Range('K3').Value = 6314
$ws.Range('D4').Value = 1976
$ws.Range('F4').Value = 1916
$ws.Range('K4').Value = 1325
$ws.Range('K5').Value = 450
$ws.Range('K6').Value = 6960
$ws.Range('D7').Value = 28166
$ws.Range('F7').Value = 24109
$ws.Range('K7').Value = 21185

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 185
$ws.Range('K6').Value = 146
$ws.Range('K7').Value = 619
$ws.Range('K8').Value = 1391
$ws.Range('K19').Value = 616
$ws.Range('K20').Value = 505
$ws.Range('K27').Value = 200
$ws.Range('K29').Value = 1150
$ws.Range('K33').Value = 924
$ws.Range('K36').Value = 273
$ws.Range('K37').Value = 718
$ws.Range('K42').Value = 786
$ws.Range('K44').Value = 179
$ws.Range('K48').Value = 267
$ws.Range('K49').Value = 115
$ws.Range('K51').Value = 273
$ws.Range('K52').Value = 560
$ws.Range('K54').Value = 412
$ws.Range('D63').Value = 355
$ws.Range('F63').Value = 200
$ws.Range('J63').Value = 116
$ws.Range('K63').Value = 66
$ws.Range('K64').Value = 132
$ws.Range('K65').Value = 501
$ws.Range('K67').Value = 834
$ws.Range('K71').Value = 62
$ws.Range('K72').Value = 104
$ws.Range('K76').Value = 286
$ws.Range('K79').Value = 531
$ws.Range('J83').Value = 594
$ws.Range('K83').Value = 465
$ws.Range('K84').Value = 165
$ws.Range('K85').Value = 982
$ws.Range('K86').Value = 130
$ws.Range('K88').Value = 226
$ws.Range('K89').Value = 314
$ws.Range('K90').Value = 195
$ws.Range('K94').Value = 286
$ws.Range('K96').Value = 221
$ws.Range('K97').Value = 166
$ws.Range('K100').Value = 40
$ws.Range('D101').Value = 28166
$ws.Range('F101').Value = 24109
$ws.Range('K101').Value = 21185

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 69
$ws.Range('K6').Value = 95
$ws.Range('K7').Value = 221

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K3').Value = 205
$ws.Range('K7').Value = 619

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 89
$ws.Range('K3').Value = 95
$ws.Range('K4').Value = 38
$ws.Range('K7').Value = 314

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 321
$ws.Range('K3').Value = 339
$ws.Range('K6').Value = 239
$ws.Range('K7').Value = 982

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 161
$ws.Range('K7').Value = 560

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 425
$ws.Range('K7').Value = 1391

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 160
$ws.Range('J4').Value = 26
$ws.Range('J7').Value = 594
$ws.Range('K7').Value = 465

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 239
$ws.Range('K3').Value = 332
$ws.Range('K4').Value = 47
$ws.Range('K6').Value = 282
$ws.Range('K7').Value = 924

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 239
$ws.Range('K7').Value = 718

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 168
$ws.Range('K6').Value = 178
$ws.Range('K7').Value = 501

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 230
$ws.Range('K3').Value = 300
$ws.Range('K5').Value = 20
$ws.Range('K6').Value = 238
$ws.Range('K7').Value = 834

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K6').Value = 32
$ws.Range('K7').Value = 165

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K6').Value = 58
$ws.Range('K7').Value = 115

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K4').Value = 20
$ws.Range('K7').Value = 412

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 331
$ws.Range('K3').Value = 413
$ws.Range('K4').Value = 56
$ws.Range('K6').Value = 322
$ws.Range('K7').Value = 1150

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K2').Value = 38
$ws.Range('K4').Value = 38
$ws.Range('K7').Value = 267

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 182
$ws.Range('K3').Value = 187
$ws.Range('K7').Value = 616

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K6').Value = 71
$ws.Range('K7').Value = 179

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K6').Value = 148
$ws.Range('K7').Value = 286

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K2').Value = 55
$ws.Range('K7').Value = 146

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 216
$ws.Range('K6').Value = 293
$ws.Range('K7').Value = 786

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 176
$ws.Range('K6').Value = 134
$ws.Range('K7').Value = 531

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 50
$ws.Range('K7').Value = 132

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 171
$ws.Range('K7').Value = 505

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 107
$ws.Range('K7').Value = 273

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('K3').Value = 6
$ws.Range('K7').Value = 40

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K6').Value = 126
$ws.Range('K7').Value = 286

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K3').Value = 48
$ws.Range('K7').Value = 185

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K6').Value = 92
$ws.Range('K7').Value = 166

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K6').Value = 94
$ws.Range('K7').Value = 226

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K3').Value = 46
$ws.Range('K4').Value = 27
$ws.Range('K7').Value = 200

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K6').Value = 32
$ws.Range('K7').Value = 130

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 73
$ws.Range('K7').Value = 195

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 73
$ws.Range('K6').Value = 88
$ws.Range('K7').Value = 273

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K3').Value = 21
$ws.Range('K7').Value = 62

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K4').Value = 8
$ws.Range('K7').Value = 104
